# Add the new "2022-Q3" quarterly sheet (sheetId 2) by cloning the
# structure/formatting of the existing "2022-Q2" sheet, inserting it
# immediately before "2022-Q2", then overwriting its data with the new
# quarter's holdings and trimming away the extra (stale) rows.

$wb = $excel.ActiveWorkbook

$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The cloned sheet still has 2022-Q2's 18 data rows (rows 2-18); the new
# quarter only has 5 holdings, so clear out everything below row 6.
$q3.Range("A7:H18").Clear()

$q3data = @(
  @(0, "004814", "中欧红利优享灵活配置混合A", "17.00", "89.59", "3.41", "0.5797", 6),
  @(1, "004815", "中欧红利优享灵活配置混合C", "8.56",  "89.59", "3.41", "0.2919", 6),
  @(2, "515150", "富国中证国企一带一路ETF",   "5.83",  "99.44", "2.08", "0.1213", 5),
  @(3, "515110", "易方达中证国企一带一路ETF", "3.60",  "99.20", "2.01", "0.0724", 6),
  @(4, "515990", "汇添富中证国企一带一路ETF", "0.93",  "98.74", "2.01", "0.0187", 6)
)

$r = 2
foreach ($row in $q3data) {
  $q3.Range("A$r").Value = $row[0]
  $q3.Range("B$r").Value = "'" + $row[1]
  $q3.Range("C$r").Value = "'" + $row[2]
  $q3.Range("D$r").Value = "'" + $row[3]
  $q3.Range("E$r").Value = "'" + $row[4]
  $q3.Range("F$r").Value = "'" + $row[5]
  $q3.Range("G$r").Value = "'" + $row[6]
  $q3.Range("H$r").Value = $row[7]
  $r++
}

# Prepend a "2022-Q3" row to the "总计" (summary) sheet and renumber the
# existing index column (A) as everything shifts down one row.
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()
$total.Range("A2").Value = 0
$total.Range("A2").Style = $total.Range("A3").Style
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 1.08

$row = 3
$idx = 1
while ($row -le 9) {
  $total.Range("A$row").Value = $idx
  $row++
  $idx++
}
